$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.109.84'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '2.306.05'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.05'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.519'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.12%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.516'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.93'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0791'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.01'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.01%  '
$ws.Range('E14').Value = '  -1.44%  '
$ws.Range('D15').Value = '2.668.16'
$ws.Range('E15').Value = '  +0.14%  '
$ws.Range('D16').Value = '2.280.33'
$ws.Range('E16').Value = '  -2.83%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.789'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').Value = '43.053.94'
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.16'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.80%  '
$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('E21').Value = '  -1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.40'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.74%  '
$ws.Range('E24').Value = '  -2.42%  '
$ws.Range('E25').Value = '  -0.41%  '
$ws.Range('E26').Value = '  -1.59%  '
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '166.96'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('E30').Value = '  -11.81%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.98'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.69%  '
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.83%  '
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0689'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('E39').Value = '  -0.68%  '
$ws.Range('E40').Value = '  +1.14%  '
$ws.Range('E41').Value = '  -3.30%  '
$ws.Range('D42').Value = '2.012.30'
$ws.Range('E42').Value = '  +0.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0287'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('E44').Value = '  -7.12%  '
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.43'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.38%  '
$ws.Range('E47').Value = '  -2.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '2.538.65'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.58%  '
